$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly-documented alternate-function / timer mappings for the
# Arduino-header pin rows (column D), and for the LED3/LED4 + LED2 rows.
$ws.Range("D44").Value = "T3C4,T8C4"
$ws.Range("D45").Value = "T1C2N,T8C2N,T15C1,I2C2SDA,SPI2MISO"

$ws.Range("D48").Value = "USART3RX"
$ws.Range("D49").Value = "USART3TX"

$ws.Range("D54").Value = "TIM2C2,TIM5C2,TIM15C1N,UART4RX"
$ws.Range("D55").Value = "TIM2C1,TIM5C1,UART4TX"
$ws.Range("D56").Value = "T4C3"
$ws.Range("D57").Value = "TIM1C2,TIM3C3,TIM8C2N"
$ws.Range("D58").Value = "TIM2C4,TIM5C4,TIM15C2,UART2RX"
$ws.Range("D59").Value = "T3C1,SPI1MISO,SPI3MISO"
$ws.Range("D60").Value = "T1C3N,T3C4,T8C3N"
$ws.Range("D61").Value = "USART2CK,SPI1NSS,SPI3NSS"
$ws.Range("D62").Value = "I2C3_SMBA"
$ws.Range("D63").Value = "T2C1,SPI1NSS,SPI3NSS"
$ws.Range("D64").Value = "T2C3,T5C3,T15C1,UART2TX"
$ws.Range("D65").Value = "T1C1N,T3C2,T8C1N,T17C1,SPI1MOSI"
$ws.Range("D66").Value = "T3C1,T16C1,SPI1MISO"
$ws.Range("D67").Value = "T2C1,T8C1,SPI1SCK"
$ws.Range("D68").Value = "T4C4,T17C1,CAN1TX,I2C1SDA,SPI2NSS"
$ws.Range("D69").Value = "T4C3,T16C1,I2C1SCL"

# Clear out the stray PUSHBUTTON/PB13 leftover row at the bottom - A70 is
# removed entirely, C70 keeps its border style but loses its value.
$ws.Range("A70").Clear()
$ws.Range("C70").ClearContents()

# Match the author's final cursor position/selection.
$ws.Range("C44").Select()
